$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/number-like cell updates -------------------------------
$ws.Range("E2").Value = "2026-02-16 19:48:36"
$ws.Range("I2").Value = "20.4 mm"
$ws.Range("E3").Value = "2026-02-16 19:48:38"
$ws.Range("O3").Value = "-0.9 °C"
$ws.Range("E4").Value = "2026-02-16 19:48:40"
$ws.Range("O4").Value = "13.9 °C"
$ws.Range("E5").Value = "2026-02-16 19:48:43"
$ws.Range("I5").Value = "24.4 mm"
$ws.Range("N5").Value = "-1.8 °C 19:25 TU"
$ws.Range("E6").Value = "2026-02-16 19:48:45"
$ws.Range("E7").Value = "2026-02-16 19:48:48"
$ws.Range("E8").Value = "2026-02-16 19:48:50"
$ws.Range("E9").Value = "2026-02-16 19:48:53"
$ws.Range("E10").Value = "2026-02-16 19:48:55"
$ws.Range("E11").Value = "2026-02-16 19:48:56"
$ws.Range("O11").Value = "6.8 °C"
$ws.Range("E12").Value = "2026-02-16 19:48:57"
$ws.Range("O12").Value = "10.9 °C"
$ws.Range("E13").Value = "2026-02-16 19:48:58"
$ws.Range("E14").Value = "2026-02-16 19:49:00"
$ws.Range("E15").Value = "2026-02-16 19:49:01"
$ws.Range("O15").Value = "11.5 °C"
$ws.Range("E16").Value = "2026-02-16 19:49:02"
$ws.Range("N16").Value = "-1.6 °C 19:27 TU"
$ws.Range("O16").Value = "-0.1 °C"
$ws.Range("E17").Value = "2026-02-16 19:49:03"
$ws.Range("E18").Value = "2026-02-16 19:49:04"
$ws.Range("O18").Value = "10.9 °C"
$ws.Range("E19").Value = "2026-02-16 19:49:05"
$ws.Range("E20").Value = "2026-02-16 19:49:06"
$ws.Range("I20").Value = "0.4 mm"
$ws.Range("E21").Value = "2026-02-16 19:49:07"
$ws.Range("O21").Value = "8.5 °C"
$ws.Range("E22").Value = "2026-02-16 19:49:10"
$ws.Range("E23").Value = "2026-02-16 19:49:12"
$ws.Range("I23").Value = "13.9 mm"
$ws.Range("N23").Value = "-1.9 °C 19:25 TU"
$ws.Range("E24").Value = "2026-02-16 19:49:15"
$ws.Range("E25").Value = "2026-02-16 19:49:17"
$ws.Range("I25").Value = "6.0 mm"
$ws.Range("N25").Value = "-0.6 °C 19:29 TU"
$ws.Range("E26").Value = "2026-02-16 19:49:20"
$ws.Range("E27").Value = "2026-02-16 19:49:22"
$ws.Range("E28").Value = "2026-02-16 19:49:24"
$ws.Range("J28").Value = "1012.8 hPa"
$ws.Range("O28").Value = "9.7 °C"
$ws.Range("E29").Value = "2026-02-16 19:49:27"
$ws.Range("E30").Value = "2026-02-16 19:49:29"
$ws.Range("K30").Value = "12.2 MJ/m2"
$ws.Range("E31").Value = "2026-02-16 19:49:32"
$ws.Range("E32").Value = "2026-02-16 19:49:34"
$ws.Range("E33").Value = "2026-02-16 19:49:37"
$ws.Range("E34").Value = "2026-02-16 19:49:40"
$ws.Range("L34").Value = "67.3 km/h - 151º 19:13 TU"
$ws.Range("O34").Value = "3.5 °C"
$ws.Range("E35").Value = "2026-02-16 19:49:42"
$ws.Range("J35").Value = "1016.7 hPa"
$ws.Range("O35").Value = "9.6 °C"
$ws.Range("E36").Value = "2026-02-16 19:49:45"
$ws.Range("L36").Value = "50.0 km/h - 327º 19:19 TU"
$ws.Range("O36").Value = "11.8 °C"
$ws.Range("E37").Value = "2026-02-16 19:49:47"
$ws.Range("E38").Value = "2026-02-16 19:49:49"
$ws.Range("E39").Value = "2026-02-16 19:49:52"
$ws.Range("I39").Value = "3.6 mm"
$ws.Range("N39").Value = "-1.7 °C 19:07 TU"
$ws.Range("E40").Value = "2026-02-16 19:49:55"
$ws.Range("E41").Value = "2026-02-16 19:49:57"
$ws.Range("O41").Value = "17.6 °C"
$ws.Range("E42").Value = "2026-02-16 19:49:59"
$ws.Range("E43").Value = "2026-02-16 19:50:02"
$ws.Range("O43").Value = "8.7 °C"
$ws.Range("E44").Value = "2026-02-16 19:50:04"
$ws.Range("I44").Value = "9.3 mm"
$ws.Range("E45").Value = "2026-02-16 19:50:07"
$ws.Range("I45").Value = "17.6 mm"
$ws.Range("E46").Value = "2026-02-16 19:50:10"
$ws.Range("J46").Value = "1016.9 hPa"

# --- Percentage-look cells: force text so Excel does not coerce them to --
# --- a numeric percentage (which would change the cell style/numFmt).  --
# --- Strategy: set via NumberFormat "@" + Value, then repair the style  --
# --- by pasting the (unchanged) format from a same-column donor cell.   --
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "71%"
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "76%"
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "67%"
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "74%"
$ws.Range("H18").PasteSpecial(-4122) | Out-Null
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "83%"
$ws.Range("H40").PasteSpecial(-4122) | Out-Null
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "79%"
$ws.Range("H43").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

